# Add the new Course_Section row (course_section_id 1395954 / PSY180 /
# AN130704.0807.5W) with its start/end/withdraw dates, then move the
# active selection/tab over to that sheet the way the author left the
# workbook.

$wb = $excel.ActiveWorkbook

$wsCS = $wb.Worksheets.Item("Course_Section")

$wsCS.Range("A2").Value = 1395954
$wsCS.Range("B2").Value = "PSY180"
$wsCS.Range("C2").Value = "AN130704.0807.5W"
$wsCS.Range("D2").Value = 41459
$wsCS.Range("E2").Value = 41493
$wsCS.Range("F2").Value = 41466

$wsCS.Range("A2:C2").HorizontalAlignment = -4131
$wsCS.Range("D2:F2").HorizontalAlignment = -4131
$wsCS.Range("D2:F2").NumberFormat = "MM/DD/YY"

# Selections as left by the author: Person scrolled/selected at B2,
# Course_Section becomes the active (selected) tab with C6 selected.
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("B2").Select()

$wsCS.Activate()
$wsCS.Range("C6").Select()
